# Updated cryptos list with GitHub Actions: refresh Price (column D) and
# Volume(1h) (column E) values for each coin row, and apply the Maker /
# TrustWalletToken row swap (rows 40-41) with their refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.207.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.681.68'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5251'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06350'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07643'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.690.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E12').Style = 'Normal'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5746'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008237'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.236.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.007'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.865'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.240'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.008'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.774'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1263'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06286'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.61%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.378'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.314'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.566'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.570'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.681'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.023'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6113'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.421'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.746'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.180'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01610'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Maker'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.098.68'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8862'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.832.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.035'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05269'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4281'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.008'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.25%  '
$ws.Range('E51').Style = 'Normal'

Write-Host "Crypto list updated"
